$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest crypto data.
# Values are forced to remain text (matching the original inlineStr cells)
# by using a leading apostrophe, then the style is reset to Normal so no
# stray 'quote prefix' formatting is left on the cell.

$ws.Range("D2").Value = "'55.033.66"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.04%  "
$ws.Range("D3").Value = "'2.276.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'506.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.54%  "
$ws.Range("D6").Value = "'128.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = "  -0.41%  "
$ws.Range("E8").Value = "  +0.96%  "
$ws.Range("D9").Value = "'2.288.68"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.23%  "
$ws.Range("E10").Value = "  +4.21%  "
$ws.Range("E11").Value = "  +1.00%  "
$ws.Range("E12").Value = "  +7.05%  "
$ws.Range("E13").Value = "  +2.47%  "
$ws.Range("D14").Value = "'23.68"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.88%  "
$ws.Range("D15").Value = "'2.682.53"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.71%  "
$ws.Range("D16").Value = "'55.066.56"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.14%  "
$ws.Range("E17").Value = "  +1.57%  "
$ws.Range("D18").Value = "'2.267.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.50%  "
$ws.Range("D19").Value = "'10.37"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.50%  "
$ws.Range("E20").Value = "  +1.61%  "
$ws.Range("D21").Value = "'314.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.28%  "
$ws.Range("D22").Value = "'6.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.19%  "
$ws.Range("E23").Value = "  -0.32%  "
$ws.Range("E24").Value = "  -1.63%  "
$ws.Range("E25").Value = "  -0.49%  "
$ws.Range("E26").Value = "  +4.88%  "
$ws.Range("D27").Value = "'7.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.15%  "
$ws.Range("D28").Value = "'171.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.75%  "
$ws.Range("E29").Value = "  +4.30%  "
$ws.Range("E30").Value = "  +2.84%  "
$ws.Range("E31").Value = "  +2.92%  "
$ws.Range("E32").Value = "  +7.63%  "
$ws.Range("D34").Value = "'17.97"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.47%  "
$ws.Range("D35").Value = "'0.993"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.55%  "
$ws.Range("E36").Value = "  +3.94%  "
$ws.Range("D37").Value = "'0.903"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.94%  "
$ws.Range("D38").Value = "'3.89"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.34%  "
$ws.Range("D39").Value = "'36.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.71%  "
$ws.Range("E40").Value = "  +4.93%  "
$ws.Range("E41").Value = "  +1.17%  "
$ws.Range("D42").Value = "'136.73"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +9.75%  "
$ws.Range("E43").Value = "  +3.85%  "
$ws.Range("E44").Value = "  +2.27%  "
$ws.Range("D45").Value = "'259.23"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.19%  "
$ws.Range("E47").Value = "  +3.76%  "
$ws.Range("E48").Value = "  +1.51%  "
$ws.Range("E49").Value = "  +4.68%  "
$ws.Range("E50").Value = "  +1.35%  "
$ws.Range("D51").Value = "'16.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.97%  "
